$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New PSSM score values (rows 2-21, columns B-K) from the updated supplemental figures run
$data = @{
  2 = @(-19.63064929472359, 2.34690769624868, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  3 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.410467388885449, -19.63064929472359, -19.63064929472359)
  4 = @(-19.63064929472359, 2.317244785837842, 2.775532090706277, -19.63064929472359, 4.321926406453217, -19.63064929472359, 1.794915456780632, -19.63064929472359, 2.395040802752221, -19.63064929472359)
  5 = @(-19.63064929472359, 0.8947424726078854, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.108892774100241, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  6 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  7 = @(2.912870877285962, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  8 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, 2.799318489513292, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  9 = @(3.640253519789197, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  10 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 1.470190747489741, -19.63064929472359, 2.217098177207498)
  11 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, 2.063982775497589, -19.63064929472359, 2.514609029372532, -19.63064929472359, -19.63064929472359, -19.63064929472359, 1.423963634149038)
  12 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  13 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, 1.734915446696921, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.293753313878272, 1.491193409966272)
  14 = @(-19.63064929472359, -19.63064929472359, 1.806079019719414, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.027881679386507)
  15 = @(-19.63064929472359, -19.63064929472359, -1.061828473312973, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359)
  16 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.313688517040358, -19.63064929472359)
  17 = @(-19.63064929472359, 0.8023135776802848, 0.2274364518266718, -19.63064929472359, -19.63064929472359, -19.63064929472359, 0.4974268624678322, 0.9197038893019943, 1.217090914443931, -19.63064929472359)
  18 = @(-19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, -19.63064929472359, 0.4198118553492086, 0.8912587659986965, 1.345090239868354, -19.63064929472359)
  19 = @(-19.63064929472359, -19.63064929472359, 1.673426240783866, -19.63064929472359, -19.63064929472359, -19.63064929472359, 1.824170578699149, 2.100106487992451, -19.63064929472359, -19.63064929472359)
  20 = @(-19.63064929472359, 1.609688445899752, 2.267925892700912, -19.63064929472359, -19.63064929472359, -19.63064929472359, 2.220427550413943, 1.95568520427186, -19.63064929472359, 2.530550569823757)
  21 = @(-19.63064929472359, 1.710978651661286, -19.63064929472359, 2.466985842868191, -19.63064929472359, 3.31787567880144, 2.480067535408778, -19.63064929472359, -19.63064929472359, -19.63064929472359)
}

foreach ($row in $data.Keys) {
  $values = $data[$row]
  for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2  # column B is index 2
    $ws.Cells.Item($row, $col).Value = $values[$i]
  }
}
